# Saldo_guide.xlsx edit: refresh the daily extraction -
# rename the sheet/tab to the new run timestamp, bump every
# balance date (column G) from 2024-08-01 to 2024-08-02, and
# roll the "pending" (D) amount into the "settled" (E) amount
# for the rows that cleared since the last extraction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet/tab to the new extraction run id.
$ws.Name = "IClientBalance-20240802-093319-"

# Every balance row's date (col G) moves from 45505 (2024-08-01)
# to 45506 (2024-08-02).
$ws.Range("G2:G274").Value = 45506

# Rows whose pending balance (D) fully cleared into the settled
# balance (E) since the prior extraction: D resets to 0, and E
# becomes the (unchanged) total already shown in H.
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 4708.26

$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0

$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 705.89

$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 2468.17

$ws.Range("D99").Value = 0
$ws.Range("E99").Value = -2917.45

$ws.Range("D104").Value = 0
$ws.Range("E104").Value = 17515

$ws.Range("D108").Value = 0
$ws.Range("E108").Value = -2732.93

$ws.Range("D132").Value = 0
$ws.Range("E132").Value = -2242.4899999999998

$ws.Range("D173").Value = 0
$ws.Range("E173").Value = -2961.52

$ws.Range("D235").Value = 0
$ws.Range("E235").Value = 3697.14

$ws.Range("D249").Value = 0
$ws.Range("E249").Value = -4768.5

$ws.Range("D264").Value = 0
$ws.Range("E264").Value = -1176.96

$ws.Range("D265").Value = 0
$ws.Range("E265").Value = -1939.72

$ws.Range("D273").Value = 0
$ws.Range("E273").Value = -1052.67

# Restore the view: no frozen top-left scroll, single cell
# selection on N18 (matches the author's last saved cursor spot).
$ws.Range("N18").Select()
